$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Hot Potato" sheet: append rows 42-44 (new scrim results)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Hot Potato")

# Row 42 - style pattern copied from an existing "Equipo 2" row (row 4)
$ws.Range("A4:N4").Copy()
$ws.Range("A42:N42").PasteSpecial(-4122)
$ws.Range("A42").Value2 = "OTIS"
$ws.Range("B42").Value2 = "RUFFS"
$ws.Range("C42").Value2 = "HANK"
$ws.Range("D42").Value2 = "KAZE"
$ws.Range("E42").Value2 = "CROW"
$ws.Range("F42").Value2 = "GRIFF"
$ws.Range("G42").Value2 = "Equipo 2"
$ws.Range("H42").Value2 = "TRB|Zeus 解開"
$ws.Range("I42").Value2 = "TRB|Lxffy"
$ws.Range("J42").Value2 = "TRB|R B M"
$ws.Range("K42").Value2 = "NHG|Xemp"
$ws.Range("L42").Value2 = "KCP|Fade"
$ws.Range("M42").Value2 = "KCP|Tyrant"
$ws.Range("N42").Value2 = "20250723T210916.000Z"

# Row 43 - style pattern copied from an existing "Equipo 1" row (row 6)
$ws.Range("A6:N6").Copy()
$ws.Range("A43:N43").PasteSpecial(-4122)
$ws.Range("A43").Value2 = "EMZ"
$ws.Range("B43").Value2 = "CHUCK"
$ws.Range("C43").Value2 = "RICO"
$ws.Range("D43").Value2 = "KAZE"
$ws.Range("E43").Value2 = "CHARLIE"
$ws.Range("F43").Value2 = "BULL"
$ws.Range("G43").Value2 = "Equipo 1"
$ws.Range("H43").Value2 = "TRB|Lxffy"
$ws.Range("I43").Value2 = "TRB|R B M"
$ws.Range("J43").Value2 = "TRB|Zeus 解開"
$ws.Range("K43").Value2 = "KCP|Fade"
$ws.Range("L43").Value2 = "KCP|Zoulan"
$ws.Range("M43").Value2 = "KCP|Tyrant"
$ws.Range("N43").Value2 = "20250723T205905.000Z"

# Row 44 - same "Equipo 1" style
$ws.Range("A6:N6").Copy()
$ws.Range("A44:N44").PasteSpecial(-4122)
$ws.Range("A44").Value2 = "EMZ"
$ws.Range("B44").Value2 = "CHUCK"
$ws.Range("C44").Value2 = "RICO"
$ws.Range("D44").Value2 = "KAZE"
$ws.Range("E44").Value2 = "CHARLIE"
$ws.Range("F44").Value2 = "BULL"
$ws.Range("G44").Value2 = "Equipo 1"
$ws.Range("H44").Value2 = "TRB|Lxffy"
$ws.Range("I44").Value2 = "TRB|R B M"
$ws.Range("J44").Value2 = "TRB|Zeus 解開"
$ws.Range("K44").Value2 = "KCP|Fade"
$ws.Range("L44").Value2 = "KCP|Zoulan"
$ws.Range("M44").Value2 = "KCP|Tyrant"
$ws.Range("N44").Value2 = "20250723T205709.000Z"

# ---------------------------------------------------------------------------
# "Layer Cake" sheet: append rows 44-47 (new scrim results)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Layer Cake")

# Row 44 - style pattern copied from an existing "Equipo 1" row (row 8)
$ws2.Range("A8:N8").Copy()
$ws2.Range("A44:N44").PasteSpecial(-4122)
$ws2.Range("A44").Value2 = "MR. P"
$ws2.Range("B44").Value2 = "FINX"
$ws2.Range("C44").Value2 = "EMZ"
$ws2.Range("D44").Value2 = "DOUG"
$ws2.Range("E44").Value2 = "CARL"
$ws2.Range("F44").Value2 = "KAZE"
$ws2.Range("G44").Value2 = "Equipo 1"
$ws2.Range("H44").Value2 = "TRB|Zeus 解開"
$ws2.Range("I44").Value2 = "TRB|R B M"
$ws2.Range("J44").Value2 = "TRB|Lxffy"
$ws2.Range("K44").Value2 = "KCP|Zoulan"
$ws2.Range("L44").Value2 = "KCP|Tyrant"
$ws2.Range("M44").Value2 = "KCP|Fade"
$ws2.Range("N44").Value2 = "20250723T205134.000Z"

# Row 45 - same "Equipo 1" style
$ws2.Range("A8:N8").Copy()
$ws2.Range("A45:N45").PasteSpecial(-4122)
$ws2.Range("A45").Value2 = "MR. P"
$ws2.Range("B45").Value2 = "FINX"
$ws2.Range("C45").Value2 = "EMZ"
$ws2.Range("D45").Value2 = "DOUG"
$ws2.Range("E45").Value2 = "CARL"
$ws2.Range("F45").Value2 = "KAZE"
$ws2.Range("G45").Value2 = "Equipo 1"
$ws2.Range("H45").Value2 = "TRB|Zeus 解開"
$ws2.Range("I45").Value2 = "TRB|R B M"
$ws2.Range("J45").Value2 = "TRB|Lxffy"
$ws2.Range("K45").Value2 = "KCP|Zoulan"
$ws2.Range("L45").Value2 = "KCP|Tyrant"
$ws2.Range("M45").Value2 = "KCP|Fade"
$ws2.Range("N45").Value2 = "20250723T204933.000Z"

# Row 46 - same "Equipo 1" style
$ws2.Range("A8:N8").Copy()
$ws2.Range("A46:N46").PasteSpecial(-4122)
$ws2.Range("A46").Value2 = "ALLI"
$ws2.Range("B46").Value2 = "MANDY"
$ws2.Range("C46").Value2 = "CORDELIUS"
$ws2.Range("D46").Value2 = "BUSTER"
$ws2.Range("E46").Value2 = "CHESTER"
$ws2.Range("F46").Value2 = "LOU"
$ws2.Range("G46").Value2 = "Equipo 1"
$ws2.Range("H46").Value2 = "TRB|R B M"
$ws2.Range("I46").Value2 = "TRB|Zeus 解開"
$ws2.Range("J46").Value2 = "TRB|Lxffy"
$ws2.Range("K46").Value2 = "KCP|Fade"
$ws2.Range("L46").Value2 = "KCP|Zoulan"
$ws2.Range("M46").Value2 = "KCP|Tyrant"
$ws2.Range("N46").Value2 = "20250723T204346.000Z"

# Row 47 - same "Equipo 1" style
$ws2.Range("A8:N8").Copy()
$ws2.Range("A47:N47").PasteSpecial(-4122)
$ws2.Range("A47").Value2 = "ALLI"
$ws2.Range("B47").Value2 = "MANDY"
$ws2.Range("C47").Value2 = "CORDELIUS"
$ws2.Range("D47").Value2 = "BUSTER"
$ws2.Range("E47").Value2 = "CHESTER"
$ws2.Range("F47").Value2 = "LOU"
$ws2.Range("G47").Value2 = "Equipo 1"
$ws2.Range("H47").Value2 = "TRB|R B M"
$ws2.Range("I47").Value2 = "TRB|Zeus 解開"
$ws2.Range("J47").Value2 = "TRB|Lxffy"
$ws2.Range("K47").Value2 = "KCP|Fade"
$ws2.Range("L47").Value2 = "KCP|Zoulan"
$ws2.Range("M47").Value2 = "KCP|Tyrant"
$ws2.Range("N47").Value2 = "20250723T204126.000Z"

$excel.CutCopyMode = $false
